$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.58%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.067"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.20%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08062"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.89%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.926"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.141"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.22%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.839"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.72%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9297"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.05%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1281"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.85%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1909"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.74%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09202"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.24%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03477"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.11%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09907"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.65%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001414"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.73%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006630"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "13.91%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.610"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.22%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.77%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3422"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.34%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1336"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.41%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.164"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.35%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2531"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.32%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04409"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.19%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004722"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.83%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.58%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003130"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.40%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05145"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.24%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007630"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.42%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01010"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.84%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1366"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.58%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01070"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.31%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006304"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.70%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.96"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.46%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001601"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.54%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"
